$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.992.93"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "2.236.33"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "'305.77"
$ws.Range("E5").Value = "  -4.23%  "
$ws.Range("D6").Value = "'94.68"
$ws.Range("E6").Value = "  -5.99%  "
$ws.Range("D7").Value = "'0.569"
$ws.Range("E7").Value = "  -0.90%  "
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("E9").Value = "  -4.93%  "
$ws.Range("D10").Value = "'34.72"
$ws.Range("E10").Value = "  -5.90%  "
$ws.Range("D11").Value = "'0.0805"
$ws.Range("E11").Value = "  -2.90%  "
$ws.Range("D12").Value = "'7.20"
$ws.Range("E12").Value = "  -4.43%  "
$ws.Range("E13").Value = "  -1.12%  "
$ws.Range("D14").Value = "2.575.55"
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("D15").Value = "2.235.42"
$ws.Range("E15").Value = "  -0.65%  "
$ws.Range("D16").Value = "'0.824"
$ws.Range("E16").Value = "  -3.29%  "
$ws.Range("D17").Value = "'13.62"
$ws.Range("E17").Value = "  -4.94%  "
$ws.Range("D18").Value = "43.864.24"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").Value = "0.0₃0957"
$ws.Range("E19").Value = "  -2.02%  "
$ws.Range("D20").Value = "'12.08"
$ws.Range("E20").Value = "  -10.37%  "
$ws.Range("D21").Value = "'6.24"
$ws.Range("E21").Value = "  -3.22%  "
$ws.Range("D22").Value = "'64.90"
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("D23").Value = "'236.33"
$ws.Range("E23").Value = "  +1.14%  "
$ws.Range("D24").Value = "'2.93"
$ws.Range("E24").Value = "  -5.28%  "
$ws.Range("E25").Value = "  -5.21%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").Value = "'9.98"
$ws.Range("E27").Value = "  -6.11%  "
$ws.Range("D28").Value = "'37.50"
$ws.Range("E28").Value = "  -3.30%  "
$ws.Range("E29").Value = "  -1.10%  "
$ws.Range("D30").Value = "'5.92"
$ws.Range("E30").Value = "  -3.41%  "
$ws.Range("D31").Value = "'19.84"
$ws.Range("E31").Value = "  -1.38%  "
$ws.Range("D32").Value = "'153.06"
$ws.Range("E32").Value = "  -4.10%  "
$ws.Range("D33").Value = "'0.0800"
$ws.Range("E33").Value = "  -5.23%  "
$ws.Range("D34").Value = "'3.22"
$ws.Range("E34").Value = "  +5.48%  "
$ws.Range("D35").Value = "'2.57"
$ws.Range("E35").Value = "  -3.75%  "
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("E37").Value = "  -6.69%  "
$ws.Range("D38").Value = "'1.77"
$ws.Range("E38").Value = "  -8.95%  "
$ws.Range("D39").Value = "'15.11"
$ws.Range("E39").Value = "  -8.25%  "
$ws.Range("D40").Value = "'3.81"
$ws.Range("E40").Value = "  -8.49%  "
$ws.Range("E41").Value = "  -9.60%  "
$ws.Range("E42").Value = "  -4.28%  "
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("D44").Value = "1.726.46"
$ws.Range("E44").Value = "  -2.04%  "
$ws.Range("D45").Value = "'85.27"
$ws.Range("E45").Value = "  +5.15%  "
$ws.Range("E46").Value = "  -4.28%  "
$ws.Range("D47").Value = "'99.74"
$ws.Range("E47").Value = "  -3.79%  "
$ws.Range("E48").Value = "  -5.12%  "
$ws.Range("D49").Value = "'8.05"
$ws.Range("E49").Value = "  -2.86%  "
$ws.Range("D50").Value = "'68.64"
$ws.Range("E50").Value = "  -8.00%  "
$ws.Range("D51").Value = "'54.09"
$ws.Range("E51").Value = "  -5.51%  "
